$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.646.64"
$ws.Range("E2").Value = "  -3.07%  "
$ws.Range("D3").Value = "2.095.78"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.32"
$ws.Range("E5").Value = "  -2.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5123"
$ws.Range("E7").Value = "  -2.44%  "
$ws.Range("E8").Value = "  -2.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.53"
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09129"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.83"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").Value = "2.096.43"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.739"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.193"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.66"
$ws.Range("E16").Value = "  -2.25%  "
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.008"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.10"
$ws.Range("E19").Value = "  +8.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06645"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.175"
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("D23").Value = "29.700.28"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.57"
$ws.Range("E24").Value = "  -1.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.308"
$ws.Range("E25").Value = "  -3.13%  "
$ws.Range("D26").Value = "2.346.12"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("E27").Value = "  -2.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.35"
$ws.Range("E28").Value = "  -1.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.521"
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("E30").Value = "  -3.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.129"
$ws.Range("E31").Value = "  -5.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1043"
$ws.Range("E32").Value = "  -3.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.635"
$ws.Range("E33").Value = "  -1.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.160"
$ws.Range("E34").Value = "  -3.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.962"
$ws.Range("E35").Value = "  -1.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.047"
$ws.Range("E36").Value = "  +2.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.27"
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06664"
$ws.Range("E39").Value = "  -2.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2230"
$ws.Range("E40").Value = "  -3.55%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6847"
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.36"
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.290"
$ws.Range("E43").Value = "  +1.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6685"
$ws.Range("E44").Value = "  +3.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.11"
$ws.Range("E45").Value = "  -4.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.291"
$ws.Range("E46").Value = "  -1.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.606"
$ws.Range("E47").Value = "  -4.02%  "
$ws.Range("E48").Value = "  -2.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "81.76"
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("E50").Value = "  -8.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.162"
$ws.Range("E51").Value = "  -2.54%  "
